$d = $word.ActiveDocument

# --- Paragraph 1: "Apple carga en su iPhone el libro electrónico de Amazon"
#     -> "Cuando los Astros se Alinean" + bookmark "_GoBack" at the very end
#     of the paragraph (collapsed, before the paragraph mark).
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$r1.End = $r1.End - 1
# Insert a temporary trailing marker character so the bookmark insertion
# point isn't exactly at the paragraph's text end (avoids an edge case),
# then remove the marker once the bookmark is in place.
$r1.Text = "Cuando los Astros se AlineanX"

$bmPos = $r1.Start + 28
$rb = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $rb)

$delRange = $d.Range($bmPos, $bmPos + 1)
$delRange.Delete()

# --- Paragraph 2
$p2 = $d.Paragraphs(2)
$r2 = $p2.Range
$r2.End = $r2.End - 1
$r2.Text = "Los eclipses tienen lugar cuando un cuerpo celeste queda oculto parcial o totalmente debido a otro."

# --- Paragraph 3 (originally split across runs with proofErr tags around "Touch")
$p3 = $d.Paragraphs(3)
$r3 = $p3.Range
$r3.End = $r3.End - 1
$r3.Text = "Son fenómenos relativamente comunes en el caso del planeta Tierra: los más conocidos popularmente son los eclipses solares y lunares. Los eclipses de la Luna son bastante frecuentes. Un eclipse lunar sucede porque la sombra de la Tierra se proyecta sobre la Luna."

# --- Paragraph 4 (originally split across runs with proofErr tags around "The")
$p4 = $d.Paragraphs(4)
$r4 = $p4.Range
$r4.End = $r4.End - 1
$r4.Text = "Los eclipses solares se deben a la interposición de la Luna entre la Tierra y el Sol."

# --- Paragraph 5 (originally split across many runs with proofErr tags)
$p5 = $d.Paragraphs(5)
$r5 = $p5.Range
$r5.End = $r5.End - 1
$r5.Text = "También se pueden interponer entre los dos astros los planetas Mercurio y Venus. Pero en este caso, dado que el tamaño del disco aparente de estos planetas es diminuto respecto al del Sol, se suele hablar de tránsitos, y no de eclipses."

# --- Paragraph 6 (originally multiple runs plus the "_GoBack" bookmark in the middle)
$p6 = $d.Paragraphs(6)
$r6 = $p6.Range
$r6.End = $r6.End - 1
$r6.Text = "Cuando se interponen entre la Tierra y el Sol, en realidad Mercurio y Venus no eclipsan nada."
